# Receipt_Images_Processing.pptx - "finished classification & clustering"
#
# 1) Slide 1 title: "...and Classification" -> "...and Item Classification"
# 2) Slide 2 subtitle bullets: trim to shorter phrasing (goals achieved)
# 3) Date placeholder (datetimeFigureOut field) on the slide master and
#    every slide layout: 12/16/2024 -> 12/17/2024
# 4) (best effort) empty presentation-level guide list extension

$p = $ppt.ActivePresentation

# --- 1) Title slide -------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Receipt Data Extraction and Item Classification"

# --- 2) Goals slide ---------------------------------------------------------
$s2 = $p.Slides.Item(2)
$goalsShape = $s2.Shapes.Item(2)
$goalsShape.TextFrame.TextRange.Text = "Data Extraction`rReceipt Items Classification `rFinding similar Items"

# --- 3) Date placeholders on master + every layout -------------------------
function Set-DatePlaceholderText {
    param($shapes, $newText)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "12/17/2024"

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    Set-DatePlaceholderText $cl.Shapes "12/17/2024"
}

# --- 4) Empty presentation-level slide guide list (best effort) ------------
try {
    $guides = $p.Guides
    $null = $guides.Add(1, 3000)
} catch {
}
